$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1 "Mes" - reuse the same formatting (style) as D1 (bold, bordered, centered)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Mes"

# Data cells E2:E151 "Abr-24" - plain/default formatting like the other data columns
for ($r = 2; $r -le 151; $r++) {
    $ws.Cells.Item($r, 5).Value = "Abr-24"
}
